$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 08:52"

# Countries data refreshed; a handful of countries changed rank, shifting the
# surrounding rows down by one position. Rewrite the affected rows (47-51,
# 88-90, 112-113, 128-130) plus the standalone data correction on row 179.

$ws.Cells.Item(47, 1).Value = "Ucrania"
$ws.Cells.Item(47, 2).Value = 2777
$ws.Cells.Item(47, 3).Value = 266
$ws.Cells.Item(47, 4).Value = 89
$ws.Cells.Item(47, 5).Value = 2605
$ws.Cells.Item(47, 6).Value = 45
$ws.Cells.Item(47, 7).Value = 10
$ws.Cells.Item(47, 8).Value = 83

$ws.Cells.Item(48, 1).Value = "Republica Dominicana"
$ws.Cells.Item(48, 2).Value = 2759
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 108
$ws.Cells.Item(48, 5).Value = 2516
$ws.Cells.Item(48, 6).Value = 147
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 135

$ws.Cells.Item(49, 1).Value = "Catar"
$ws.Cells.Item(49, 2).Value = 2728
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 247
$ws.Cells.Item(49, 5).Value = 2475
$ws.Cells.Item(49, 6).Value = 37
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 6

$ws.Cells.Item(50, 1).Value = "Colombia"
$ws.Cells.Item(50, 2).Value = 2709
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 214
$ws.Cells.Item(50, 5).Value = 2395
$ws.Cells.Item(50, 6).Value = 92
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 100

$ws.Cells.Item(51, 1).Value = "Tailandia"
$ws.Cells.Item(51, 2).Value = 2551
$ws.Cells.Item(51, 3).Value = 33
$ws.Cells.Item(51, 4).Value = 1218
$ws.Cells.Item(51, 5).Value = 1295
$ws.Cells.Item(51, 6).Value = 61
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 38

$ws.Cells.Item(88, 1).Value = "Oman"
$ws.Cells.Item(88, 2).Value = 599
$ws.Cells.Item(88, 3).Value = 53
$ws.Cells.Item(88, 4).Value = 109
$ws.Cells.Item(88, 5).Value = 487
$ws.Cells.Item(88, 6).Value = 3
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 3

$ws.Cells.Item(89, 1).Value = "Costa Rica"
$ws.Cells.Item(89, 2).Value = 577
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 49
$ws.Cells.Item(89, 5).Value = 525
$ws.Cells.Item(89, 6).Value = 13
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 3

$ws.Cells.Item(90, 1).Value = "Afganistan"
$ws.Cells.Item(90, 2).Value = 555
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 32
$ws.Cells.Item(90, 5).Value = 505
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 18

$ws.Cells.Item(112, 1).Value = "Georgia"
$ws.Cells.Item(112, 2).Value = 252
$ws.Cells.Item(112, 3).Value = 10
$ws.Cells.Item(112, 4).Value = 60
$ws.Cells.Item(112, 5).Value = 189
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 3

$ws.Cells.Item(113, 1).Value = "Guinea"
$ws.Cells.Item(113, 2).Value = 250
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 17
$ws.Cells.Item(113, 5).Value = 233
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 0

$ws.Cells.Item(128, 1).Value = "El Salvador"
$ws.Cells.Item(128, 2).Value = 125
$ws.Cells.Item(128, 3).Value = 7
$ws.Cells.Item(128, 4).Value = 21
$ws.Cells.Item(128, 5).Value = 98
$ws.Cells.Item(128, 6).Value = 3
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 6

$ws.Cells.Item(129, 1).Value = "Camboya"
$ws.Cells.Item(129, 2).Value = 122
$ws.Cells.Item(129, 3).Value = 2
$ws.Cells.Item(129, 4).Value = 77
$ws.Cells.Item(129, 5).Value = 45
$ws.Cells.Item(129, 6).Value = 1
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 0

$ws.Cells.Item(130, 1).Value = "Ruanda"
$ws.Cells.Item(130, 2).Value = 120
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 18
$ws.Cells.Item(130, 5).Value = 102
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 0

$ws.Cells.Item(179, 1).Value = "Santa Lucia"
$ws.Cells.Item(179, 2).Value = 15
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 4
$ws.Cells.Item(179, 5).Value = 11
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 0
